# SALES fixed + grops filter  at 09/09/2022
# Append the new Tyres sales rows (116-134) to the "Holidays 2019" sheet,
# which is where the Tyres data actually lives (columns E:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$contragent = "БНХ ПОЛЬСКА"
$salesValue = 201
$dateOfSales = 44701

$rows = @(
    @{ Row=116; E="315/80R22.5"; F="BEL-158M"; G="камневыт, груз, сер" },
    @{ Row=117; E="315/80R22.5"; F="BEL-158M"; G="камневыт, груз, трп" },
    @{ Row=118; E="315/80R22.5"; F="BEL-278";  G="груз, сер" },
    @{ Row=119; E="315/80R22.5"; F="BEL-278";  G="груз, трп" },
    @{ Row=120; E="315/80R22.5"; F="BEL-268";  G="груз, сер" },
    @{ Row=121; E="315/80R22.5"; F="BEL-268";  G="груз, трп" },
    @{ Row=122; E="315/80R22.5"; F="BEL-398";  G="груз, сер" },
    @{ Row=123; E="315/80R22.5"; F="BEL-326";  G="груз, сер" },
    @{ Row=124; E="315/80R22.5"; F="BEL-326";  G="груз, трп" },
    @{ Row=125; E="315/80R22.5"; F="BEL-498";  G="156L, груз, сер" },
    @{ Row=126; E="315/80R22.5"; F="BEL-518";  G="груз, сер" },
    @{ Row=127; E="12.00R20";    F="ИД-304М";  G="16, груз, сер" },
    @{ Row=128; E="12.00R20";    F="ИД-304М";  G="18, груз, сер" },
    @{ Row=129; E="12.00R20";    F="ИД-304М";  G="16, груз, трп" },
    @{ Row=130; E="12.00R20";    F="ИД-304М";  G="18, груз, трп" },
    @{ Row=131; E="12.00R20";    F="БИ-368М";  G="18, груз, сер" },
    @{ Row=132; E="12.00R20";    F="БИ-368М";  G="18, груз, сер" },
    @{ Row=133; E="12.00R20";    F="БИ-368М";  G="18, груз, трп" },
    @{ Row=134; E="195/65R15";   F="BEL-337";  G="б/к, легк, сер" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 5).Value2 = $r.E
    $ws.Cells.Item($row, 6).Value2 = $r.F
    $ws.Cells.Item($row, 7).Value2 = $r.G
    $ws.Cells.Item($row, 8).Value2 = $salesValue
    $dateCell = $ws.Cells.Item($row, 9)
    $dateCell.Value2 = $dateOfSales
    $dateCell.NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($row, 10).Value2 = $contragent
}
